# live_trading_results.xlsx - Trade #21 closed at 2026-02-17 13:18:19 - unknown UNKNOWN +0.000%
#
# Trade #21 (MarketMaking strategy) transitions from OPEN -> CLOSED.
# Update the "All Trades" and "MarketMaking" trade logs (row 22), then
# propagate the resulting capital / P&L / trade-count figures into the
# "Strategy Status" (MarketMaking row) and "Summary" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Trade log row 22 - present identically on both "All Trades" and
#    "MarketMaking" sheets.
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G22").Value = 0.68              # Exit Price
    $ws.Range("H22").Value = "CLOSED"          # Status: OPEN -> CLOSED
    $ws.Range("I22").Value = -2.8571           # P&L %
    $ws.Range("J22").Value = -0.02             # P&L $
    $ws.Range("K22").Value = 99.26000000000001 # Capital After
    $ws.Range("P22").Value = "early_exit"      # Exit Reason
    $ws.Range("Q22").Value = 0.13              # Duration (min)
}

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.26000000000001  # Capital
$status.Range("D4").Value = 21                 # Trades
$status.Range("E4").Value = -0.75              # P&L $
$status.Range("F4").Value = -0.74              # P&L %
$status.Range("G4").Value = 33.33              # Win Rate %

# ---------------------------------------------------------------------
# 3) Summary sheet - overall roll-up figures
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.26  # Current Capital
$summary.Range("B4").Value = -0.75    # Total P&L $
$summary.Range("B5").Value = -0.71    # Total P&L %
$summary.Range("B6").Value = 21       # Total Trades
$summary.Range("B8").Value = 13       # Losing Trades
$summary.Range("B9").Value = 33.33    # Win Rate %
